# Append the newest run's profit-allocation row (10/21/2025) to the bottom
# of the data table, mirroring the existing rows 2-49.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the date column to be stored as literal text (matching the rest of
# the "Date" column, which is plain text like "10/20/2025", not a real
# date value) by briefly marking the cell as Text before assigning it -
# otherwise Excel auto-converts the date-like string into a date serial.
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = "10/21/2025"
# Reset the cell style back to the workbook default so no stray
# number-format style lingers on the cell (keeps it identical in shape to
# the untouched rows above it).
$ws.Range("A50").Style = "Normal"

$ws.Range("B50").Value = 0.1928588791428576
$ws.Range("C50").Value = 0.8071411208571424
